# Updates cryptos list with latest prices / 1h volume percentages.
# Row 39 is a new entry (MXToken); existing rows 39-50 shift down to
# 40-51 and the previous last row (Elrond) falls off the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "27.324.35", 0),
    @("E2", "  +1.72%  ", 0),
    @("D3", "1.864.08", 0),
    @("E3", "  +1.30%  ", 0),
    @("D4", "1.022", 1),
    @("E4", "  +1.39%  ", 0),
    @("D5", "313.66", 1),
    @("E5", "  +1.42%  ", 0),
    @("D6", "1.020", 1),
    @("E6", "  +1.34%  ", 0),
    @("D7", "0.4806", 1),
    @("E7", "  +1.92%  ", 0),
    @("D8", "0.3732", 1),
    @("E8", "  +2.30%  ", 0),
    @("D9", "0.07438", 1),
    @("E9", "  +4.06%  ", 0),
    @("D10", "0.9378", 1),
    @("E10", "  +2.07%  ", 0),
    @("D11", "20.72", 1),
    @("E11", "  +6.24%  ", 0),
    @("D12", "0.07891", 1),
    @("E12", "  +3.69%  ", 0),
    @("D13", "1.866.97", 0),
    @("E13", "  -4.80%  ", 0),
    @("D14", "5.439", 1),
    @("E14", "  +3.04%  ", 0),
    @("D15", "6.551", 1),
    @("E15", "  +2.47%  ", 0),
    @("D16", "90.36", 1),
    @("E16", "  +3.00%  ", 0),
    @("E17", "  +1.25%  ", 0),
    @("D18", "0.000008802", 1),
    @("E18", "  +2.09%  ", 0),
    @("D19", "1.020", 1),
    @("E19", "  +1.31%  ", 0),
    @("D20", "14.83", 1),
    @("E20", "  +2.51%  ", 0),
    @("D21", "27.371.79", 0),
    @("E21", "  +1.76%  ", 0),
    @("D22", "5.136", 1),
    @("E22", "  +2.47%  ", 0),
    @("D23", "10.72", 1),
    @("E23", "  +1.10%  ", 0),
    @("E24", "  +1.90%  ", 0),
    @("D25", "154.18", 1),
    @("E25", "  +1.76%  ", 0),
    @("D26", "18.57", 1),
    @("E26", "  +2.16%  ", 0),
    @("D27", "2.009", 1),
    @("E27", "  +0.14%  ", 0),
    @("D28", "116.04", 1),
    @("E28", "  +1.68%  ", 0),
    @("D29", "4.997", 1),
    @("E29", "  +2.89%  ", 0),
    @("D30", "0.08935", 1),
    @("E30", "  +1.32%  ", 0),
    @("E31", "  +4.11%  ", 0),
    @("D32", "1.197", 1),
    @("E32", "  +2.46%  ", 0),
    @("D33", "4.572", 1),
    @("E33", "  +2.26%  ", 0),
    @("D34", "0.7444", 1),
    @("E34", "  +0.27%  ", 0),
    @("D35", "2.687", 1),
    @("E35", "  -2.03%  ", 0),
    @("D36", "0.02054", 1),
    @("E36", "  +5.84%  ", 0),
    @("E37", "  +3.39%  ", 0),
    @("D38", "0.05302", 1),
    @("E38", "  +1.41%  ", 0),
    @("B39", "MXToken", 0),
    @("C39", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", 0),
    @("D39", "3.006", 1),
    @("E39", "  +1.34%  ", 0),
    @("B40", "TheSandbox", 0),
    @("C40", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", 0),
    @("D40", "0.5389", 1),
    @("E40", "  +4.25%  ", 0),
    @("B41", "FraxShare", 0),
    @("C41", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", 0),
    @("D41", "7.124", 1),
    @("E41", "  +2.57%  ", 0),
    @("B42", "Algorand", 0),
    @("C42", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", 0),
    @("D42", "0.1537", 1),
    @("E42", "  +1.82%  ", 0),
    @("B43", "Aptos", 0),
    @("C43", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", 0),
    @("D43", "8.403", 1),
    @("E43", "  +3.22%  ", 0),
    @("B44", "EnergySwap", 0),
    @("C44", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", 0),
    @("D44", "10.70", 1),
    @("E44", "  +2.29%  ", 0),
    @("B45", "Decentraland", 0),
    @("C45", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", 0),
    @("D45", "0.4836", 1),
    @("E45", "  +3.01%  ", 0),
    @("B46", "PaxDollar", 0),
    @("C46", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", 0),
    @("D46", "1.021", 1),
    @("E46", "  +1.47%  ", 0),
    @("B47", "NEARProtocol", 0),
    @("C47", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", 0),
    @("D47", "1.679", 1),
    @("E47", "  +5.43%  ", 0),
    @("B48", "Quant", 0),
    @("C48", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", 0),
    @("D48", "103.46", 1),
    @("E48", "  +1.73%  ", 0),
    @("B49", "Aave", 0),
    @("C49", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", 0),
    @("D49", "66.81", 1),
    @("E49", "  +2.67%  ", 0),
    @("B50", "Cronos", 0),
    @("C50", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", 0),
    @("D50", "0.06099", 1),
    @("E50", "  +1.15%  ", 0),
    @("B51", "EOS", 0),
    @("C51", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos", 0),
    @("D51", "0.9006", 1),
    @("E51", "  +1.86%  ", 0)
)

foreach ($u in $updates) {
    $ref = $u[0]
    $val = $u[1]
    $forceText = $u[2]
    $c = $ws.Range($ref)
    if ($forceText -eq 1) {
        $c.NumberFormat = "@"
        $c.Value = $val
        $c.ClearFormats()
    } else {
        $c.Value = $val
    }
}
